$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Español (column D) grades for rows 5-11 ---
$ws.Range("D5").Value = 6.5
$ws.Range("D6").Value = 7.8
$ws.Range("D7").Value = 7.7
$ws.Range("D8").Value = 8.2
$ws.Range("D9").Value = 7.5
$ws.Range("D10").Value = 6.5
$ws.Range("D11").Value = 8.7

# --- Apply "0.0" number format with centered alignment to the grade table ---
$grades = $ws.Range("B2:F11")
$grades.NumberFormat = "0.0"
$grades.HorizontalAlignment = -4108  # xlCenter

# --- Update sheet view (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("I11").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9        # xlPaperA4
$ws.PageSetup.Orientation = 1      # xlPortrait
